# Korrektur des Entwurfes von UC-Intern
# Umbenennen von Finalize in BuildGraph in der Klasse UseCase.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Method name cell: Finalize() : void -> BuildGraph() : void
$ws.Range("A41").Value = "BuildGraph() : void"

# Method description cell (B42), rewritten to reference BuildGraph instead of Finalize
$ws.Range("B42").Value = "Signales the end of the description of the use-case and builds the internal graph-representation. Before calling BuildGraph only the setters and methods starting with ""Set"" and ""Add"" are allowed to be called.`nAfter using these methods to describe the use-case, BuildGraph is called and the given description is used to build a graph representation. After calling BuildGraph only the getters are allowed to be called."

# Class description cell (B2), rewritten to reference BuildGraph instead of Finalize and tidy blank lines
$ws.Range("B2").Value = "This class is used to create and retrieve the representation of the use-case as a graph. It also provides all the information contained in a use-case.`nThe usage of the class involves two steps:`nThe first step is to build the graph. For that the setter methods and the methods starting with ""Set"" and ""Add"" are used to insert the description of a graph. For that all ""Set"" methods must be used and all setters and ""Add"" methods can be used to describe a use-case.`nThe setters (in this case setters are all methods starting with ""set_"" are not mandatory because they are simple string setters and therefor return null if no value was set.`nAfter inserting the whole description the method BuildGraph is called. After that no usage of setters, ""Set"" and ""Add"" methods is allowed anymore.`nThe second step involves retrieval of the data. For that various getters are provided with most being simple string getters. The remaining getters are used to retrive the graph in a representation with a node list and two edge matrices."

# Reset the view selection back to A1 (diff drops the stored B3 selection)
$ws.Range("A1").Select()
